$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F50").Value = 3
$ws.Range("G50").Value = 1887.18
$ws.Range("B52").Value = 6916.36
$ws.Range("F62").Value = 37
$ws.Range("G62").Value = 2624.04
$ws.Range("F68").Value = 57
$ws.Range("G68").Value = 6561.84
$ws.Range("F70").Value = 32
$ws.Range("G70").Value = 4318.4
$ws.Range("F71").Value = 365
$ws.Range("G71").Value = 23250.5
$ws.Range("F84").Value = 46
$ws.Range("G84").Value = 4713.16
$ws.Range("F85").Value = 147
$ws.Range("G85").Value = 19811.19
$ws.Range("B90").Value = 198334.97
$ws.Range("F115").Value = 228
$ws.Range("G115").Value = 22072.68
$ws.Range("B117").Value = 16124.96
$ws.Range("F140").Value = 12
$ws.Range("G140").Value = 511.68
$ws.Range("B142").Value = 3188.18
$ws.Range("F144").Value = 1194
$ws.Range("G144").Value = 10089.3
$ws.Range("F145").Value = 614
$ws.Range("G145").Value = 4905.86
$ws.Range("B147").Value = 17857.62
$ws.Range("F149").Value = 248
$ws.Range("G149").Value = 16070.4
$ws.Range("B156").Value = 34956.26
$ws.Range("F203").Value = 69
$ws.Range("G203").Value = 1391.04
$ws.Range("F214").Value = 49
$ws.Range("G214").Value = 4297.3
$ws.Range("B216").Value = 47654.28
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 123
$ws.Range("G232").Value = 5859.72
$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12
$ws.Range("F247").Value = 155
$ws.Range("G247").Value = 16106.05
$ws.Range("F249").Value = 144
$ws.Range("G249").Value = 19846.08
$ws.Range("F251").Value = 7
$ws.Range("G251").Value = 1713.25
$ws.Range("F255").Value = 600
$ws.Range("G255").Value = 102798
$ws.Range("B260").Value = 205799.3
$ws.Range("F283").Value = 42
$ws.Range("G283").Value = 14341.74
$ws.Range("F288").Value = 47
$ws.Range("G288").Value = 4370.53
$ws.Range("F291").Value = 121
$ws.Range("G291").Value = 5204.21
$ws.Range("F294").Value = 45
$ws.Range("G294").Value = 3211.2
$ws.Range("F295").Value = 5
$ws.Range("G295").Value = 518.45
$ws.Range("F296").Value = 82
$ws.Range("G296").Value = 1738.4
$ws.Range("B304").Value = 192046.05
$ws.Range("F320").Value = 68
$ws.Range("G320").Value = 4668.2
$ws.Range("F327").Value = 16
$ws.Range("G327").Value = 4024
$ws.Range("B330").Value = 31104.94
$ws.Range("F335").Value = 0
$ws.Range("G335").Value = 0
$ws.Range("F345").Value = 74
$ws.Range("G345").Value = 4544.34
$ws.Range("B346").Value = 28178.87
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29
$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32
$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68
$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945
$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65
$ws.Range("F423").Value = 5
$ws.Range("G423").Value = 551.35
$ws.Range("B424").Value = 4115.4
$ws.Range("F434").Value = 36
$ws.Range("G434").Value = 1175.04
$ws.Range("B435").Value = 1385.48
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("F455").Value = 49
$ws.Range("G455").Value = 3116.89
$ws.Range("F458").Value = 27
$ws.Range("G458").Value = 489.78
$ws.Range("B460").Value = 14774.17
$ws.Range("F474").Value = 108
$ws.Range("G474").Value = 3545.64
$ws.Range("B475").Value = 46032.89
$ws.Range("F477").Value = 15
$ws.Range("G477").Value = 680.1
$ws.Range("B478").Value = 680.1
$ws.Range("F525").Value = 7
$ws.Range("G525").Value = 637.5599999999999
$ws.Range("B526").Value = 637.5599999999999
$ws.Range("F539").Value = 48
$ws.Range("G539").Value = 12435.36
$ws.Range("B547").Value = 24323.32
$ws.Range("F551").Value = 3
$ws.Range("G551").Value = 429.39
$ws.Range("B560").Value = 8080.04
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 27
$ws.Range("G573").Value = 1103.49
$ws.Range("F575").Value = 0
$ws.Range("G575").Value = 0
$ws.Range("F577").Value = 76
$ws.Range("G577").Value = 3267.24
$ws.Range("F580").Value = 67
$ws.Range("G580").Value = 3818.33
$ws.Range("B583").Value = 26570.86
$ws.Range("F599").Value = 1973
$ws.Range("G599").Value = 321816.03
$ws.Range("F601").Value = 460
$ws.Range("G601").Value = 130120.2
$ws.Range("F602").Value = 349
$ws.Range("G602").Value = 50482.85
$ws.Range("B606").Value = 503267.13
$ws.Range("F614").Value = 0
$ws.Range("G614").Value = 0
$ws.Range("B618").Value = 46271.01
$ws.Range("B619").Value = 1944492.24
$ws.Range("B620").Value = 1944492.24
